$wb = $excel.ActiveWorkbook

# --- Sheet "Merchant(P)CustomFee_QPayNow": fix Fee Type / DisplayLabel values ---
# (was incorrectly copied from the Surcharge sheet: "Surcharge"/"SurchargeFee")
$sheetCustomFee = $wb.Worksheets.Item("Merchant(P)CustomFee_QPayNow")
$sheetCustomFee.Activate() | Out-Null
$sheetCustomFee.Range("I2").Value = "CustomFee"
$sheetCustomFee.Range("M2").Value = "CustomFee"
$sheetCustomFee.Range("H1").Select() | Out-Null

# --- Sheet "QPayWithSendlinkCustomFee": reset selection/scroll back to the top ---
$sheetSendLink = $wb.Worksheets.Item("QPayWithSendlinkCustomFee")
$sheetSendLink.Activate() | Out-Null
$sheetSendLink.Range("A1").Select() | Out-Null

# --- Sheet "RunManager": fix capitalization bug, update selection ---
# "CustomFeeTxWithQPayNow_MerchantPrimaryLogin" -> "customFeeTxWithQPayNow_MerchantPrimaryLogin"
$sheetRunManager = $wb.Worksheets.Item("RunManager")
$sheetRunManager.Activate() | Out-Null
$sheetRunManager.Range("A18").Value = "customFeeTxWithQPayNow_MerchantPrimaryLogin"
$sheetRunManager.Range("C18").Select() | Out-Null
